$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk price / volume(1h) updates for rows 2-48 ---
# NumberFormat is forced to text ("@") before each write so Excel's COM
# layer does not auto-coerce numeric-looking strings (e.g. "308.94",
# "1.002", scientific-looking small decimals) into real numbers - the
# source data must stay literal text, matching the inlineStr cells in
# the original workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.917.89"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.814.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4649"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8695"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.27"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.840.93"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.364"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07103"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.503"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008703"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.945.02"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.032.74"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.896"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.35"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.121"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.252"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.46"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08907"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7559"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.159"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.481"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.901"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.084"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05284"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.971"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.249"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5308"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.319"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1653"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.418"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4872"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.14%  "

# --- Rows 49-51: Quant/NEARProtocol swap position, Cronos values update ---
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.661"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06298"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.03%  "
